# Powerpoint writer: consolidate text run nodes.
#
# Each affected paragraph currently holds its text split across three
# runs, e.g. "The" / " " / "Moon" (the leading word, a lone space, and
# the trailing word). We merge the first two runs ("<word>" + " ")
# into a single run ("<word> "), leaving the final word as its own
# run - exactly mirroring the target OOXML (2 runs instead of 3).
#
# Selecting the exact character span that covers the first two runs
# (via TextRange.Characters(start, length)) and re-assigning its .Text
# forces the writer to rebuild just that span as one run, instead of
# doing a minimal char-level diff that would otherwise preserve the
# existing run boundaries untouched.

$p = $ppt.ActivePresentation

# --- Slide 2: TextBox "The" + " " + "Moon" -> "The " + "Moon" ---
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tb2.Characters(1, 4).Text = "The "

# --- Slide 3: Title "One" + " " + "More" -> "One " + "More" ---
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Characters(1, 4).Text = "One "

# --- Slide 3: TextBox "The" + " " + "Moon" -> "The " + "Moon" ---
$tb3 = $s3.Shapes.Item(3).TextFrame.TextRange
$tb3.Characters(1, 4).Text = "The "
